$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = "stimuli/img_r10cu.png"
$ws.Range("M2").Value = 78.52380952380952
$ws.Range("N2").Value = 56.14285714285715
$ws.Range("O2").Value = 67.33333333333333
$ws.Range("P2").Value = 42
$ws.Range("Q2").Value = 7
$ws.Range("R2").Value = 7
$ws.Range("S2").Value = 7
$ws.Range("T2").Value = 7
$ws.Range("U2").Value = 7
$ws.Range("V2").Value = 6

# Row 3
$ws.Range("L3").Value = "stimuli/img_pey7u.png"
$ws.Range("M3").Value = 30.34883720930232
$ws.Range("N3").Value = 20.34883720930232
$ws.Range("O3").Value = 25.34883720930232
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 2
$ws.Range("T3").Value = 2
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 2

# Row 4
$ws.Range("I4").ClearContents()
$ws.Range("J4").Value = "new"
$ws.Range("K4").Value = "f"
$ws.Range("L4").Value = "stimuli/img_5tr4v.png"
$ws.Range("M4").Value = 56.86046511627907
$ws.Range("N4").Value = 39.3953488372093
$ws.Range("O4").Value = 48.12790697674419
$ws.Range("P4").Value = 43
$ws.Range("Q4").Value = 4
$ws.Range("R4").Value = 4
$ws.Range("S4").Value = 4
$ws.Range("T4").Value = 4
$ws.Range("U4").Value = 4
$ws.Range("V4").Value = 4

# Row 5
$ws.Range("I5").Value = "target"
$ws.Range("J5").Value = "old"
$ws.Range("K5").Value = "j"
$ws.Range("L5").Value = "stimuli/img_njhlh.png"
$ws.Range("M5").Value = 59.74418604651163
$ws.Range("N5").Value = 41.51162790697674
$ws.Range("O5").Value = 50.62790697674419
$ws.Range("Q5").Value = 4
$ws.Range("R5").Value = 4
$ws.Range("S5").Value = 4
$ws.Range("T5").Value = 4
$ws.Range("U5").Value = 4
$ws.Range("V5").Value = 4

# Row 6
$ws.Range("I6").ClearContents()
$ws.Range("J6").Value = "new"
$ws.Range("K6").Value = "f"
$ws.Range("L6").Value = "stimuli/img_5nlnv.png"
$ws.Range("M6").Value = 86.1219512195122
$ws.Range("N6").Value = 69.1951219512195
$ws.Range("O6").Value = 77.65853658536585
$ws.Range("P6").Value = 41
$ws.Range("Q6").Value = 9
$ws.Range("R6").Value = 9
$ws.Range("S6").Value = 9
$ws.Range("T6").Value = 9
$ws.Range("U6").Value = 9
$ws.Range("V6").Value = 9

# Row 7
$ws.Range("I7").Value = "target"
$ws.Range("J7").Value = "old"
$ws.Range("K7").Value = "j"
$ws.Range("L7").Value = "stimuli/img_xbtev.png"
$ws.Range("M7").Value = 13.68181818181818
$ws.Range("N7").Value = 8.568181818181818
$ws.Range("O7").Value = 11.125
$ws.Range("P7").Value = 44
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = 1
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 1
$ws.Range("U7").Value = 1
$ws.Range("V7").Value = 1

# Row 8
$ws.Range("I8").Value = "target"
$ws.Range("J8").Value = "old"
$ws.Range("K8").Value = "j"
$ws.Range("L8").Value = "stimuli/img_3ze38.png"
$ws.Range("M8").Value = 35.53191489361702
$ws.Range("N8").Value = 28.4468085106383
$ws.Range("O8").Value = 31.98936170212766
$ws.Range("P8").Value = 47
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 2
$ws.Range("S8").Value = 2
$ws.Range("T8").Value = 3
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 3

# Row 9
$ws.Range("L9").Value = "stimuli/img_9oofc.png"
$ws.Range("M9").Value = 82.47619047619048
$ws.Range("N9").Value = 65.5
$ws.Range("O9").Value = 73.98809523809524
$ws.Range("P9").Value = 42
$ws.Range("Q9").Value = 8
$ws.Range("R9").Value = 8
$ws.Range("S9").Value = 8
$ws.Range("T9").Value = 8
$ws.Range("U9").Value = 8
$ws.Range("V9").Value = 8

# Row 10
$ws.Range("I10").ClearContents()
$ws.Range("J10").Value = "new"
$ws.Range("K10").Value = "f"
$ws.Range("L10").Value = "stimuli/img_il020.png"
$ws.Range("M10").Value = 18.85416666666667
$ws.Range("N10").Value = 16.16666666666667
$ws.Range("O10").Value = 17.51041666666667
$ws.Range("P10").Value = 48
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 1
$ws.Range("S10").Value = 1
$ws.Range("T10").Value = 1
$ws.Range("U10").Value = 1
$ws.Range("V10").Value = 1

# Row 11
$ws.Range("I11").Value = "target"
$ws.Range("J11").Value = "old"
$ws.Range("K11").Value = "j"
$ws.Range("L11").Value = "stimuli/img_mdh76.png"
$ws.Range("M11").Value = 37.31914893617022
$ws.Range("N11").Value = 25.12765957446809
$ws.Range("O11").Value = 31.22340425531915
$ws.Range("P11").Value = 47
$ws.Range("Q11").Value = 2
$ws.Range("R11").Value = 2
$ws.Range("S11").Value = 2
$ws.Range("T11").Value = 2
$ws.Range("U11").Value = 3
$ws.Range("V11").Value = 2

# Row 12
$ws.Range("I12").ClearContents()
$ws.Range("J12").Value = "new"
$ws.Range("K12").Value = "f"
$ws.Range("L12").Value = "stimuli/img_x9w7o.png"
$ws.Range("M12").Value = 92.38888888888889
$ws.Range("N12").Value = 72.94444444444444
$ws.Range("O12").Value = 82.66666666666666
$ws.Range("P12").Value = 36
$ws.Range("Q12").Value = 10
$ws.Range("R12").Value = 10
$ws.Range("S12").Value = 10
$ws.Range("T12").Value = 10
$ws.Range("U12").Value = 10
$ws.Range("V12").Value = 10

# Row 13
$ws.Range("I13").Value = "target"
$ws.Range("J13").Value = "old"
$ws.Range("K13").Value = "j"
$ws.Range("L13").Value = "stimuli/img_qdln8.png"
$ws.Range("M13").Value = 85.51162790697674
$ws.Range("N13").Value = 67.86046511627907
$ws.Range("O13").Value = 76.68604651162791
$ws.Range("P13").Value = 43

# Row 14
$ws.Range("L14").Value = "stimuli/img_2qhro.png"
$ws.Range("M14").Value = 81.73809523809524
$ws.Range("N14").Value = 62.73809523809524
$ws.Range("O14").Value = 72.23809523809524
$ws.Range("P14").Value = 42
$ws.Range("Q14").Value = 8
$ws.Range("R14").Value = 8
$ws.Range("S14").Value = 8
$ws.Range("T14").Value = 8
$ws.Range("U14").Value = 8
$ws.Range("V14").Value = 8

# Row 15
$ws.Range("L15").Value = "stimuli/img_qz292.png"
$ws.Range("M15").Value = 78.26666666666667
$ws.Range("N15").Value = 59.13333333333333
$ws.Range("O15").Value = 68.7
$ws.Range("P15").Value = 45
$ws.Range("Q15").Value = 7
$ws.Range("R15").Value = 7
$ws.Range("S15").Value = 7
$ws.Range("T15").Value = 7
$ws.Range("U15").Value = 7
$ws.Range("V15").Value = 7

# Row 16
$ws.Range("I16").Value = "target"
$ws.Range("J16").Value = "old"
$ws.Range("K16").Value = "j"
$ws.Range("L16").Value = "stimuli/img_9684y.png"
$ws.Range("M16").Value = 77.95744680851064
$ws.Range("N16").Value = 56.70212765957447
$ws.Range("O16").Value = 67.32978723404256
$ws.Range("T16").Value = 6

# Row 17
$ws.Range("I17").ClearContents()
$ws.Range("J17").Value = "new"
$ws.Range("K17").Value = "f"
$ws.Range("L17").Value = "stimuli/img_4o8l0.png"
$ws.Range("M17").Value = 46.02173913043478
$ws.Range("N17").Value = 31.45652173913043
$ws.Range("O17").Value = 38.73913043478261
$ws.Range("P17").Value = 46
$ws.Range("Q17").Value = 3
$ws.Range("R17").Value = 3
$ws.Range("S17").Value = 3
$ws.Range("T17").Value = 3
$ws.Range("U17").Value = 3
$ws.Range("V17").Value = 3

# Row 18
$ws.Range("I18").ClearContents()
$ws.Range("J18").Value = "new"
$ws.Range("K18").Value = "f"
$ws.Range("L18").Value = "stimuli/img_s2zoe.png"
$ws.Range("M18").Value = 64.71428571428571
$ws.Range("N18").Value = 44.90476190476191
$ws.Range("O18").Value = 54.80952380952381
$ws.Range("Q18").Value = 5
$ws.Range("R18").Value = 5
$ws.Range("S18").Value = 5
$ws.Range("T18").Value = 5
$ws.Range("U18").Value = 5
$ws.Range("V18").Value = 5

# Row 19
$ws.Range("I19").Value = "target"
$ws.Range("J19").Value = "old"
$ws.Range("K19").Value = "j"
$ws.Range("L19").Value = "stimuli/img_kq9s9.png"
$ws.Range("M19").Value = 62.30232558139535
$ws.Range("N19").Value = 39.97674418604651
$ws.Range("O19").Value = 51.13953488372093
$ws.Range("P19").Value = 43
$ws.Range("Q19").Value = 4
$ws.Range("R19").Value = 4
$ws.Range("S19").Value = 4
$ws.Range("T19").Value = 5
$ws.Range("U19").Value = 5
$ws.Range("V19").Value = 4

# Row 20
$ws.Range("I20").ClearContents()
$ws.Range("J20").Value = "new"
$ws.Range("K20").Value = "f"
$ws.Range("L20").Value = "stimuli/img_jpjeg.png"
$ws.Range("M20").Value = 90.90697674418605
$ws.Range("N20").Value = 74.3953488372093
$ws.Range("O20").Value = 82.65116279069767
$ws.Range("P20").Value = 43

# Row 21
$ws.Range("L21").Value = "stimuli/img_rru0v.png"
$ws.Range("M21").Value = 56.45238095238095
$ws.Range("N21").Value = 39.42857142857143
$ws.Range("O21").Value = 47.94047619047619
$ws.Range("P21").Value = 42
$ws.Range("Q21").Value = 4
$ws.Range("R21").Value = 4
$ws.Range("S21").Value = 4
$ws.Range("T21").Value = 4
$ws.Range("U21").Value = 4
$ws.Range("V21").Value = 4

# Row 22
$ws.Range("L22").Value = "stimuli/img_vg73h.png"
$ws.Range("M22").Value = 87.7
$ws.Range("N22").Value = 72.4
$ws.Range("O22").Value = 80.05000000000001
$ws.Range("P22").Value = 50
$ws.Range("Q22").Value = 10
$ws.Range("R22").Value = 10
$ws.Range("S22").Value = 10
$ws.Range("T22").Value = 10
$ws.Range("U22").Value = 10
$ws.Range("V22").Value = 10

# Row 23
$ws.Range("I23").Value = "target"
$ws.Range("J23").Value = "old"
$ws.Range("K23").Value = "j"
$ws.Range("L23").Value = "stimuli/img_99exx.png"
$ws.Range("M23").Value = 70.02272727272727
$ws.Range("N23").Value = 51.88636363636363
$ws.Range("O23").Value = 60.95454545454545
$ws.Range("T23").Value = 5
$ws.Range("U23").Value = 5

# Row 27
$ws.Range("I27").ClearContents()
$ws.Range("J27").Value = "new"
$ws.Range("K27").Value = "f"
$ws.Range("L27").Value = "stimuli/img_iudc4.png"
$ws.Range("M27").Value = 73.625
$ws.Range("N27").Value = 52.275
$ws.Range("O27").Value = 62.95
$ws.Range("P27").Value = 40
$ws.Range("Q27").Value = 6
$ws.Range("R27").Value = 6
$ws.Range("S27").Value = 6
$ws.Range("T27").Value = 6
$ws.Range("U27").Value = 6
$ws.Range("V27").Value = 6

# Row 28
$ws.Range("L28").Value = "stimuli/img_bbs77.png"
$ws.Range("M28").Value = 31.64444444444445
$ws.Range("N28").Value = 21.26666666666667
$ws.Range("O28").Value = 26.45555555555556
$ws.Range("P28").Value = 45
$ws.Range("Q28").Value = 2
$ws.Range("R28").Value = 2
$ws.Range("S28").Value = 2
$ws.Range("T28").Value = 2
$ws.Range("U28").Value = 2
$ws.Range("V28").Value = 2

# Row 29
$ws.Range("L29").Value = "stimuli/catch_08.jpg"

# Row 32
$ws.Range("L32").Value = "stimuli/img_lzz3x.png"
$ws.Range("M32").Value = 18.46341463414634
$ws.Range("N32").Value = 11.92682926829268
$ws.Range("O32").Value = 15.19512195121951
$ws.Range("P32").Value = 41
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = 1
$ws.Range("S32").Value = 1
$ws.Range("T32").Value = 1
$ws.Range("U32").Value = 1
$ws.Range("V32").Value = 1

# Row 33
$ws.Range("L33").Value = "stimuli/img_196rk.png"
$ws.Range("M33").Value = 86.53488372093024
$ws.Range("N33").Value = 69.46511627906976
$ws.Range("O33").Value = 78
$ws.Range("P33").Value = 43
$ws.Range("Q33").Value = 9
$ws.Range("R33").Value = 9
$ws.Range("S33").Value = 9
$ws.Range("T33").Value = 9
$ws.Range("U33").Value = 9
$ws.Range("V33").Value = 9

# Row 34
$ws.Range("L34").Value = "stimuli/img_37hgm.png"
$ws.Range("M34").Value = 70.95454545454545
$ws.Range("N34").Value = 54.77272727272727
$ws.Range("O34").Value = 62.86363636363636
$ws.Range("P34").Value = 44
$ws.Range("Q34").Value = 6
$ws.Range("R34").Value = 6
$ws.Range("S34").Value = 6
$ws.Range("T34").Value = 6
$ws.Range("U34").Value = 6
$ws.Range("V34").Value = 6

# Row 37
$ws.Range("L37").Value = "stimuli/img_tbs4n.png"
$ws.Range("M37").Value = 78.95744680851064
$ws.Range("N37").Value = 58.97872340425532
$ws.Range("O37").Value = 68.96808510638297
$ws.Range("P37").Value = 47
$ws.Range("Q37").Value = 7
$ws.Range("R37").Value = 7
$ws.Range("S37").Value = 7
$ws.Range("T37").Value = 7
$ws.Range("U37").Value = 7
$ws.Range("V37").Value = 7

# Row 39
$ws.Range("L39").Value = "stimuli/img_rg4in.png"
$ws.Range("M39").Value = 49.3695652173913
$ws.Range("N39").Value = 30.21739130434782
$ws.Range("O39").Value = 39.79347826086956
$ws.Range("P39").Value = 46
$ws.Range("Q39").Value = 3
$ws.Range("R39").Value = 3
$ws.Range("S39").Value = 3
$ws.Range("T39").Value = 3
$ws.Range("U39").Value = 3
$ws.Range("V39").Value = 3

# Row 41
$ws.Range("L41").Value = "stimuli/img_eiu3c.png"
$ws.Range("M41").Value = 65.1590909090909
$ws.Range("N41").Value = 46.22727272727273
$ws.Range("O41").Value = 55.69318181818181
$ws.Range("P41").Value = 44
$ws.Range("Q41").Value = 5
$ws.Range("R41").Value = 5
$ws.Range("S41").Value = 5
$ws.Range("T41").Value = 5
$ws.Range("U41").Value = 5
$ws.Range("V41").Value = 5

# Row 42
$ws.Range("L42").Value = "stimuli/img_bj99b.png"
$ws.Range("M42").Value = 82.79069767441861
$ws.Range("N42").Value = 65.46511627906976
$ws.Range("O42").Value = 74.12790697674419
$ws.Range("P42").Value = 43
$ws.Range("Q42").Value = 8
$ws.Range("R42").Value = 8
$ws.Range("S42").Value = 8
$ws.Range("T42").Value = 8
$ws.Range("U42").Value = 8
$ws.Range("V42").Value = 8
